$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a number-looking string (e.g. "602.68") must be
# force-typed as text so Excel does not silently convert them to a Number,
# then the temporary text format is cleared so the cell keeps the default
# (General) style, matching cells that already contained unambiguous text.
$textForceCells = @("D5", "D6", "D14", "D20", "D23", "D25", "D26", "D29", "D30", "D31", "D32", "D38", "D39", "D43", "D44", "D47", "D51")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.708.81"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.614.42"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "602.68"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").Value = "154.58"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").Value = "2.613.14"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("E10").Value = "  +8.27%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "28.04"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").Value = "3.088.18"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "67.564.75"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "2.619.65"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "365.95"
$ws.Range("E20").Value = "  +3.53%  "
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").Value = "2.10"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "69.97"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "10.13"
$ws.Range("E26").Value = "  -3.02%  "
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("D28").Value = "2.746.20"
$ws.Range("D29").Value = "583.81"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "7.94"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").Value = "19.40"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "155.17"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("D43").Value = "2.66"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("D44").Value = "41.12"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D47").Value = "156.56"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "0.0₆0286"
$ws.Range("E48").Value = "  -7.32%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").Value = "0.624"
$ws.Range("E51").Value = "  +0.56%  "

# Restore default (General) formatting on the cells we temporarily forced to
# text, now that the literal value is safely stored as text.
foreach ($c in $textForceCells) {
    $ws.Range($c).ClearFormats()
}
